# Update the natmi LR-pair stats (Cadm3-Cadm3) sheet following Dr Hou's advice:
# the Ligand/Receptor-expressing cell counts increased from 1 to 3, which changes
# the average expression values (G, M) and, downstream of those, the total
# expression, specificity and edge-weight columns (H..T) for every row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns E..T (5..20) for rows 2..10, in column order:
# E  Ligand-expressing cells
# F  Ligand detection rate
# G  Ligand average expression value
# H  Ligand total expression value
# I  Ligand derived specificity of average expression value
# J  Ligand derived specificity of total expression value
# K  Receptor-expressing cells
# L  Receptor detection rate
# M  Receptor average expression value
# N  Receptor total expression value
# O  Receptor derived specificity of average expression value
# P  Receptor derived specificity of total expression value
# Q  Edge average expression weight
# R  Edge total expression weight
# S  Edge average expression derived specificity
# T  Edge total expression derived specificity

$rows = @{
    2  = @(3, 1, 4.370833,         13.112499,        0.0996525630698175, 0.09965256306981748, 3, 1, 4.370833,         13.112499,        0.0996525630698175, 0.09965256306981748, 19.104181113889,   171.937630025001,  0.009930633326383954, 0.009930633326383952)
    3  = @(3, 1, 4.370833,         13.112499,        0.0996525630698175, 0.09965256306981748, 3, 1, 37.55834333333333, 112.67503,        0.8563093528905953, 0.8563093528905952,  164.1612464666633, 1477.45121819997,  0.08533342179620465,  0.08533342179620464)
    4  = @(3, 1, 4.370833,         13.112499,        0.0996525630698175, 0.09965256306981748, 3, 1, 1.931542,         5.794626,         0.04403808403958729,0.04403808403958728, 8.442447514486,    75.98202763037401, 0.004388507947228896, 0.004388507947228894)
    5  = @(3, 1, 37.55834333333333,112.67503,        0.8563093528905953, 0.8563093528905952,  3, 1, 4.370833,         13.112499,        0.0996525630698175, 0.09965256306981748, 164.1612464666633, 1477.45121819997,  0.08533342179620465,  0.08533342179620464)
    6  = @(3, 1, 37.55834333333333,112.67503,        0.8563093528905953, 0.8563093528905952,  3, 1, 37.55834333333333, 112.67503,        0.8563093528905953, 0.8563093528905952,  1410.629153944544, 12695.6623855009,  0.73326570784791,     0.7332657078479099)
    7  = @(3, 1, 37.55834333333333,112.67503,        0.8563093528905953, 0.8563093528905952,  3, 1, 1.931542,         5.794626,         0.04403808403958729,0.04403808403958728, 72.54551759875334, 652.9096583887799, 0.03771022324648064,  0.03771022324648064)
    8  = @(3, 1, 1.931542,         5.794626,         0.04403808403958729,0.04403808403958728, 3, 1, 4.370833,         13.112499,        0.0996525630698175, 0.09965256306981748, 8.442447514486,    75.98202763037401, 0.004388507947228896, 0.004388507947228894)
    9  = @(3, 1, 1.931542,         5.794626,         0.04403808403958729,0.04403808403958728, 3, 1, 37.55834333333333, 112.67503,        0.8563093528905953, 0.8563093528905952,  72.54551759875334, 652.9096583887799, 0.03771022324648064,  0.03771022324648064)
    10 = @(3, 1, 1.931542,         5.794626,         0.04403808403958729,0.04403808403958728, 3, 1, 1.931542,         5.794626,         0.04403808403958729,0.04403808403958728, 3.730854497764001, 33.577690479876,   0.001939352845877753, 0.001939352845877752)
}

foreach ($r in $rows.Keys) {
    $values = $rows[$r]
    for ($i = 0; $i -lt $values.Length; $i++) {
        $col = 5 + $i   # column E is index 5
        $ws.Cells.Item($r, $col).Value = $values[$i]
    }
}
